# [ANV] updating decay chain spreadsheet
#
# 1. Adds a new "HDPE Density" worksheet (same layout/formulas as the
#    existing "Shotcrete Density" sheet, but for polyethylene: just C & H).
# 2. Makes the new sheet the active tab (was "Shotcrete Target Fractions").
# 3. Updates stale selections left on "Shotcrete Density" /
#    "Shotcrete Target Fractions" from the prior editing session.

$wb = $excel.ActiveWorkbook

$density = $wb.Worksheets.Item("Shotcrete Density")
$targetFractions = $wb.Worksheets.Item("Shotcrete Target Fractions")

# --- Create the new sheet after the last existing sheet -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "HDPE Density"

# --- Column widths (mirror "Shotcrete Density") ----------------------------
$ws.Columns.Item(1).ColumnWidth = 20.67
$ws.Columns.Item(2).ColumnWidth = 20.67
$ws.Columns.Item(3).ColumnWidth = 20.83
$ws.Columns.Item(4).ColumnWidth = 20.67
$ws.Columns.Item(5).ColumnWidth = 20.83
$ws.Columns.Item(6).ColumnWidth = 28.5
$ws.Columns.Item(7).ColumnWidth = 21
$ws.Columns.Item(8).ColumnWidth = 23
$ws.Columns.Item(9).ColumnWidth = 31.33

# --- Header row (same headers/style as "Shotcrete Density" row 1) ---------
$density.Range("A1:I1").Copy($ws.Range("A1:I1"))

# --- Row 2: Carbon -----------------------------------------------------
$ws.Range("A2").Value = "C"
$ws.Range("B2").Value = 6
$ws.Range("C2").Formula = "=(2*G2)/(2*G2+4*G3)"
$ws.Range("D2").Value = 961
$ws.Range("E2").Formula = "=D2*(1000)*(1/100000)"
$ws.Range("F2").Formula = "=`$E`$2*(C2/100)"
$ws.Range("G2").Value = 12.011
$ws.Range("H2").Formula = "=(F2/G2)*6.0221408E+23"
$ws.Range("I2").Formula = "=H2/`$H`$4"

# --- Row 3: Hydrogen -----------------------------------------------------
$ws.Range("A3").Value = "H"
$ws.Range("B3").Value = 1
$ws.Range("C3").Formula = "=4*G3/(2*G2+4*G3)"
$ws.Range("F3").Formula = "=`$E`$2*(C3/100)"
$ws.Range("G3").Value = 1.00784
$ws.Range("H3").Formula = "=(F3/G3)*6.0221408E+23"
$ws.Range("I3").Formula = "=H3/`$H`$4"

# --- Row 4: Totals (same style as "Shotcrete Density" row 13) -------------
$density.Range("A13").Copy($ws.Range("A4"))
$density.Range("C13").Copy($ws.Range("C4"))
$density.Range("F13").Copy($ws.Range("F4"))
$density.Range("H13").Copy($ws.Range("H4"))
$density.Range("I13").Copy($ws.Range("I4"))
$ws.Range("C4").Formula = "=SUM(C2:C3)"
$ws.Range("F4").Formula = "=SUM(F2:F3)"
$ws.Range("H4").Formula = "=SUM(H2:H3)"
$ws.Range("I4").Formula = "=SUM(I2:I3)"

# --- References ------------------------------------------------------------
$ws.Range("A7").Value = "Reference: https://www.xometry.com/resources/materials/high-density-polyethylene-hdpe/"
$ws.Range("A8").Value = "https://en.wikipedia.org/wiki/High-density_polyethylene"

# --- Small summary table (same style as "Shotcrete Density" row 21) -------
$density.Range("A21:B21").Copy($ws.Range("A12:B12"))
$ws.Range("A13").Value = 6
$ws.Range("B13").Value = 0.33333333333333337
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 0.66666666666666674

# Sort the little table by Z ascending (1, 6) - matches the stored
# <sortState> left behind by the author's last "Data > Sort" in Excel.
$sortKey = $ws.Range("A13:A14")
$sortRange = $ws.Range("A13:B14")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Apply()

# --- Selections left by the prior editing session --------------------------
$density.Activate()
$density.Range("B27").Select()

$targetFractions.Activate()
$targetFractions.Range("I32").Select()

$ws.Activate()
$ws.Range("F27").Select()
